$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.069.55"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.55%  '

$ws.Range("D3").Value = "'1.763.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.07%  '

$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.38%  '

$ws.Range("D5").Value = "'335.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.91%  '

$ws.Range("D6").Value = "'0.9983"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.23%  '

$ws.Range("D7").Value = "'0.3925"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.48%  '

$ws.Range("D8").Value = "'0.3397"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.19%  '

$ws.Range("D9").Value = "'45.35"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.37%  '

$ws.Range("D10").Value = "'1.124"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.22%  '

$ws.Range("D11").Value = "'0.07236"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.97%  '

$ws.Range("B12").Value = 'Solana'
$ws.Range("C12").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D12").Value = "'22.44"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.62%  '

$ws.Range("B13").Value = 'BinanceUSD'
$ws.Range("C13").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D13").Value = "'0.9994"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.28%  '

$ws.Range("D14").Value = "'6.157"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.51%  '

$ws.Range("D15").Value = "'7.119"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.97%  '

$ws.Range("D16").Value = "'1.757.83"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.49%  '

$ws.Range("D17").Value = "'0.00001062"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.22%  '

$ws.Range("D18").Value = "'0.06613"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.94%  '

$ws.Range("D19").Value = "'80.44"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.03%  '

$ws.Range("D20").Value = "'0.9985"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.13%  '

$ws.Range("D21").Value = "'17.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.48%  '

$ws.Range("D22").Value = "'6.243"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.99%  '

$ws.Range("D23").Value = "'28.064.42"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.62%  '

$ws.Range("E24").Value = '  -3.20%  '

$ws.Range("D25").Value = "'2.392"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.54%  '

$ws.Range("D26").Value = "'154.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.17%  '

$ws.Range("D27").Value = "'19.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.39%  '

$ws.Range("D28").Value = "'2.326"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.68%  '

$ws.Range("D29").Value = "'1.943.58"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.17%  '

$ws.Range("D30").Value = "'1.284"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -10.58%  '

$ws.Range("D31").Value = "'129.55"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.87%  '

$ws.Range("D32").Value = "'4.077"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.55%  '

$ws.Range("D33").Value = "'5.837"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.13%  '

$ws.Range("D34").Value = "'0.08743"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.54%  '

$ws.Range("D35").Value = "'12.11"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.87%  '

$ws.Range("D36").Value = "'0.06209"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.69%  '

$ws.Range("D37").Value = "'0.02295"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.64%  '

$ws.Range("D38").Value = "'5.154"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.53%  '

$ws.Range("D39").Value = "'0.6518"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.89%  '

$ws.Range("D40").Value = "'0.2119"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.02%  '

$ws.Range("D41").Value = "'1.495"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.22%  '

$ws.Range("D42").Value = "'1.205"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.30%  '

$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = "'7.922"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.03%  '

$ws.Range("B44").Value = 'Frax'
$ws.Range("C44").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D44").Value = "'0.9984"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.11%  '

$ws.Range("D45").Value = "'13.84"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.53%  '

$ws.Range("D46").Value = "'3.828"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.20%  '

$ws.Range("D47").Value = "'0.6017"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.45%  '

$ws.Range("D48").Value = "'127.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.69%  '

$ws.Range("D49").Value = "'1.997"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.85%  '

$ws.Range("E50").Value = '  -4.06%  '

$ws.Range("D51").Value = "'0.07015"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.50%  '
